$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 14078
$ws1.Range("F4").Value = 680
$ws1.Range("F6").Value = 533
$ws1.Range("F7").Value = 1467

# Sheet "全部类型" (sheet4): update "想去人数" (F column) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 14078
$ws4.Range("F4").Value = 680
$ws4.Range("F8").Value = 533
$ws4.Range("F9").Value = 1467
